$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1418
$ws1.Range("F7").Value = 974
$ws1.Range("F8").Value = 1547
$ws1.Range("F10").Value = 4
$ws1.Range("F11").Value = 1455
$ws1.Range("F12").Value = 3080
$ws1.Range("F14").Value = 1765
$ws1.Range("F15").Value = 1801
$ws1.Range("F16").Value = 848
$ws1.Range("F19").Value = 1469
$ws1.Range("F23").Value = 1207
$ws1.Range("F24").Value = 401
$ws1.Range("F25").Value = 452
$ws1.Range("F26").Value = 106
$ws1.Range("F27").Value = 4770
$ws1.Range("F31").Value = 1645

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 30
$ws2.Range("F7").Value = 68

# Sheet "本地生活"
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 40

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 40
$ws4.Range("F5").Value = 30
$ws4.Range("F10").Value = 68
$ws4.Range("F12").Value = 1418
$ws4.Range("F15").Value = 974
$ws4.Range("F16").Value = 1547
$ws4.Range("F20").Value = 4
$ws4.Range("F21").Value = 1455
$ws4.Range("F22").Value = 3080
$ws4.Range("F24").Value = 1765
$ws4.Range("F25").Value = 1801
$ws4.Range("F26").Value = 848
$ws4.Range("F29").Value = 1469
$ws4.Range("F35").Value = 1207
$ws4.Range("F36").Value = 401
$ws4.Range("F37").Value = 452
$ws4.Range("F38").Value = 106
$ws4.Range("F39").Value = 4770
$ws4.Range("F43").Value = 1645
